$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update NIM, and change the execution-date cell from a date value
# to a plain text value ("25-8-2025") while keeping column C's time value.
$ws.Range("A2").Value = 2341760196
$ws.Range("B2").Value = "25-8-2025"
$ws.Range("C2").Value = 0.375

# Row 3: add a new student's schedule row, reusing row 2's cell formatting
# for the date (B) and time (C) columns.
$ws.Range("A3").Value = 2341760036
$ws.Range("B3").Value = "25-8-2025"
$ws.Range("C3").Value = 0.33333333333333331

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the saved selection to match the target workbook state.
$ws.Range("C4").Select()
